$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029170723896049
$ws.Range("D2").Value = 1.037835783583608
$ws.Range("E2").Value = 1.029037989954252
$ws.Range("F2").Value = 1.044729243049703
$ws.Range("I2").Value = 1.033768717741986
$ws.Range("J2").Value = 1.034318980332658
$ws.Range("K2").Value = 1.040625533420604
$ws.Range("L2").Value = 1.031853063021243
$ws.Range("M2").Value = 1.047499474445901
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030200869198805
$ws.Range("D3").Value = 1.038657014648209
$ws.Range("E3").Value = 1.029914439594855
$ws.Range("F3").Value = 1.045740564816937
$ws.Range("I3").Value = 1.033981050292772
$ws.Range("J3").Value = 1.034989631895561
$ws.Range("K3").Value = 1.041256623974853
$ws.Range("L3").Value = 1.032537389848686
$ws.Range("M3").Value = 1.048321571600923
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030867524096433
$ws.Range("D4").Value = 1.039188047687956
$ws.Range("E4").Value = 1.030482019661196
$ws.Range("F4").Value = 1.046395116288663
$ws.Range("I4").Value = 1.034116531359075
$ws.Range("J4").Value = 1.035423105956784
$ws.Range("K4").Value = 1.041663976785046
$ws.Range("L4").Value = 1.032980019327809
$ws.Range("M4").Value = 1.04885307395356
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031147805300519
$ws.Range("D5").Value = 1.039411207336348
$ws.Range("E5").Value = 1.030720739350374
$ws.Range("F5").Value = 1.046670327040172
$ws.Range("I5").Value = 1.034173029308594
$ws.Range("J5").Value = 1.035605222486323
$ws.Range("K5").Value = 1.041834986658392
$ws.Range("L5").Value = 1.033166058136727
$ws.Range("M5").Value = 1.049076409291361
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031194866886334
$ws.Range("D6").Value = 1.039448671717826
$ws.Range("E6").Value = 1.030760827808359
$ws.Range("F6").Value = 1.046716538304767
$ws.Range("I6").Value = 1.034182488674353
$ws.Range("J6").Value = 1.035635793838337
$ws.Range("K6").Value = 1.041863685802563
$ws.Range("L6").Value = 1.033197292355741
$ws.Range("M6").Value = 1.049113901910541
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030871269151223
$ws.Range("D7").Value = 1.039191029899622
$ws.Range("E7").Value = 1.030485209017903
$ws.Range("F7").Value = 1.046398793522436
$ws.Range("I7").Value = 1.034117288088831
$ws.Range("J7").Value = 1.035425539863059
$ws.Range("K7").Value = 1.041666262776122
$ws.Range("L7").Value = 1.032982505355666
$ws.Range("M7").Value = 1.048856058597453
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029518849257446
$ws.Range("D8").Value = 1.038113396163743
$ws.Range("E8").Value = 1.029334094821863
$ws.Range("F8").Value = 1.045070990721332
$ws.Range("I8").Value = 1.033840872184959
$ws.Range("J8").Value = 1.034545729746843
$ws.Range("K8").Value = 1.040839021124068
$ws.Range("L8").Value = 1.032084370533426
$ws.Range("M8").Value = 1.047777398709599
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027136346468486
$ws.Range("D9").Value = 1.03621175980663
$ws.Range("E9").Value = 1.027309228892544
$ws.Range("F9").Value = 1.042732476654477
$ws.Range("I9").Value = 1.033339167094681
$ws.Range("J9").Value = 1.032991718883911
$ws.Range("K9").Value = 1.039373646951323
$ws.Range("L9").Value = 1.030500420181634
$ws.Range("M9").Value = 1.04587324751237
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025548442063575
$ws.Range("D10").Value = 1.034942228615518
$ws.Range("E10").Value = 1.025961743025359
$ws.Range("F10").Value = 1.041174339981538
$ws.Range("I10").Value = 1.03299489093017
$ws.Range("J10").Value = 1.031953268537033
$ws.Range("K10").Value = 1.038391610744337
$ws.Range("L10").Value = 1.029443594723981
$ws.Range("M10").Value = 1.044601552548474
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024860964033242
$ws.Range("D11").Value = 1.034392096100501
$ws.Range("E11").Value = 1.025378850034463
$ws.Range("F11").Value = 1.040499863339695
$ws.Range("I11").Value = 1.032843494906371
$ws.Range("J11").Value = 1.031503033497974
$ws.Range("K11").Value = 1.037965170676135
$ws.Range("L11").Value = 1.028985779773875
$ws.Range("M11").Value = 1.044050364923433
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024605618235262
$ws.Range("D12").Value = 1.034187690373553
$ws.Range("E12").Value = 1.025162424940486
$ws.Range("F12").Value = 1.040249364198648
$ws.Range("I12").Value = 1.03278691115623
$ws.Range("J12").Value = 1.031335709453685
$ws.Range("K12").Value = 1.03780659023645
$ws.Range("L12").Value = 1.028815696879221
$ws.Range("M12").Value = 1.043845549306438
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024660390137375
$ws.Range("D13").Value = 1.034231538891812
$ws.Range("E13").Value = 1.025208844897214
$ws.Range("F13").Value = 1.040303095681586
$ws.Range("I13").Value = 1.032799064328614
$ws.Range("J13").Value = 1.031371604947755
$ws.Range("K13").Value = 1.037840614489968
$ws.Range("L13").Value = 1.028852181590492
$ws.Range("M13").Value = 1.043889486564429
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024839856781968
$ws.Range("D14").Value = 1.034375201120709
$ws.Range("E14").Value = 1.025360958484085
$ws.Range("F14").Value = 1.040479156357121
$ws.Range("I14").Value = 1.032838824783727
$ws.Range("J14").Value = 1.031489204207223
$ws.Range("K14").Value = 1.037952066069669
$ws.Range("L14").Value = 1.02897172127925
$ws.Range("M14").Value = 1.044033436429955
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024950433976319
$ws.Range("D15").Value = 1.034463707963043
$ws.Range("E15").Value = 1.025454692273032
$ws.Range("F15").Value = 1.040587637339215
$ws.Range("I15").Value = 1.032863276339818
$ws.Range("J15").Value = 1.031561649512383
$ws.Range("K15").Value = 1.038020711014845
$ws.Range("L15").Value = 1.029045369658178
$ws.Range("M15").Value = 1.044122118106277
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025594069708064
$ws.Range("D16").Value = 1.034978730366154
$ws.Range("E16").Value = 1.026000439926246
$ws.Range("F16").Value = 1.041219107170271
$ws.Range("I16").Value = 1.033004889690496
$ws.Range("J16").Value = 1.03198313697771
$ws.Range("K16").Value = 1.038419886673474
$ws.Range("L16").Value = 1.029473974160547
$ws.Range("M16").Value = 1.044638121814536
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025997830757015
$ws.Range("D17").Value = 1.035301679098512
$ws.Range("E17").Value = 1.026342927980669
$ws.Range("F17").Value = 1.041615267197935
$ws.Range("I17").Value = 1.03309309854415
$ws.Range("J17").Value = 1.032247370086398
$ws.Range("K17").Value = 1.038669954885548
$ws.Range("L17").Value = 1.029742772491174
$ws.Range("M17").Value = 1.04496165429739
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026233346925997
$ws.Range("D18").Value = 1.035490009339718
$ws.Range("E18").Value = 1.026542751128645
$ws.Range("F18").Value = 1.041846360565334
$ws.Range("I18").Value = 1.033144325196647
$ws.Range("J18").Value = 1.032401436884477
$ws.Range("K18").Value = 1.038815698480284
$ws.Range("L18").Value = 1.029899538440151
$ws.Range("M18").Value = 1.045150313642158
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026313653414345
$ws.Range("D19").Value = 1.035554218231125
$ws.Range("E19").Value = 1.026610895081259
$ws.Range("F19").Value = 1.041925160797834
$ws.Range("I19").Value = 1.03316175414272
$ws.Range("J19").Value = 1.032453960189376
$ws.Range("K19").Value = 1.038865373426626
$ws.Range("L19").Value = 1.029952988246194
$ws.Range("M19").Value = 1.045214632824444
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025954510068853
$ws.Range("D20").Value = 1.035267033918498
$ws.Range("E20").Value = 1.026306176492176
$ws.Range("F20").Value = 1.04157276091185
$ws.Range("I20").Value = 1.033083657743921
$ws.Range("J20").Value = 1.03221902615455
$ws.Range("K20").Value = 1.03864313702754
$ws.Range("L20").Value = 1.029713934995424
$ws.Range("M20").Value = 1.044926947665989
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024787007938587
$ws.Range("D21").Value = 1.034332897881728
$ws.Range("E21").Value = 1.025316162409044
$ws.Range("F21").Value = 1.040427310002883
$ws.Range("I21").Value = 1.032827125938902
$ws.Range("J21").Value = 1.031454576547364
$ws.Range("K21").Value = 1.03791925136696
$ws.Range("L21").Value = 1.028936520646414
$ws.Range("M21").Value = 1.043991048987115
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024053034828347
$ws.Range("D22").Value = 1.033745210990955
$ws.Range("E22").Value = 1.024694206732743
$ws.Range("F22").Value = 1.03970730210371
$ws.Range("I22").Value = 1.032663817919297
$ws.Range("J22").Value = 1.030973435247005
$ws.Range("K22").Value = 1.037463065200833
$ws.Range("L22").Value = 1.028447555791217
$ws.Range("M22").Value = 1.043402149746786
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024442120149381
$ws.Range("D23").Value = 1.034056788695649
$ws.Range("E23").Value = 1.025023869097191
$ws.Range("F23").Value = 1.040088974464161
$ws.Range("I23").Value = 1.03275058159461
$ws.Range("J23").Value = 1.031228544786316
$ws.Range("K23").Value = 1.037704997542567
$ws.Range("L23").Value = 1.028706781637326
$ws.Range("M23").Value = 1.04371438008193
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025974084794577
$ws.Range("D24").Value = 1.035282688708523
$ws.Range("E24").Value = 1.026322782735873
$ws.Range("F24").Value = 1.041591967611225
$ws.Range("I24").Value = 1.033087924327515
$ws.Range("J24").Value = 1.03223183372997
$ws.Range("K24").Value = 1.038655255224557
$ws.Range("L24").Value = 1.02972696547888
$ws.Range("M24").Value = 1.044942630258399
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027752204946177
$ws.Range("D25").Value = 1.036703694417066
$ws.Range("E25").Value = 1.027832280927728
$ws.Range("F25").Value = 1.043336887402372
$ws.Range("I25").Value = 1.033470600226364
$ws.Range("J25").Value = 1.033393900471737
$ws.Range("K25").Value = 1.039753386351722
$ws.Range("L25").Value = 1.030910062842991
$ws.Range("M25").Value = 1.046365917262673
